# Project DesignFirst: update rule R30's "Integer min" threshold (C10) to 100.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 100
